# Delete the record with Kayıt No (record id) 11225399 from both the
# "Kayitlar" master sheet and the matching "Merkez İlçe" sheet, shifting
# the remaining rows up (matches commit: "Kayıt silindi: 11225399").

$wb = $excel.ActiveWorkbook

$recordId = "11225399"
$sheetNames = @("Kayitlar", "Merkez İlçe")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $foundCell = $ws.Columns.Item(1).Find($recordId)
    if ($foundCell -ne $null) {
        $ws.Rows.Item($foundCell.Row).Delete()
    }
}
